$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 10917.3, 11579, 9364, 0.2034403324127197),
    @(1, 11066.46666666667, 11859, 9396, 0.2115513801574707),
    @(2, 10921.36666666667, 11558, 9991, 0.2185593207677205),
    @(3, 11558.53333333333, 12256, 10469, 0.2077418168385824),
    @(4, 10512.93333333333, 11319, 9611, 0.2034905513127645),
    @(5, 11466.9, 12369, 9922, 0.2134225845336914),
    @(6, 11013.23333333333, 11835, 10138, 0.2154336055119832),
    @(7, 10743.9, 11548, 10143, 0.2091625134150187),
    @(8, 10653.63333333333, 11415, 9151, 0.2059109767278035),
    @(9, 10609.66666666667, 11597, 9665, 0.2058092594146728)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
